$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4979.6665
$ws.Range("J43").Value = 6469.75
$ws.Range("L43").Value = 6469.75
$ws.Range("N43").Value = -6607.75

$ws.Range("H53").Value = 209.93333
$ws.Range("I53").Value = 200.5
$ws.Range("J53").Value = 220.71428
$ws.Range("K53").Value = 200.5
$ws.Range("L53").Value = 220.71428
$ws.Range("M53").Value = 436.5
$ws.Range("N53").Value = -1494.71428

$ws.Range("H76").Value = 7485.5
$ws.Range("I76").Value = 6699.75
$ws.Range("J76").Value = 7799.8
$ws.Range("K76").Value = 6699.75
$ws.Range("L76").Value = 7799.8
$ws.Range("M76").Value = -6384.75
$ws.Range("N76").Value = -8429.799999999999

$ws.Range("H79").Value = 7485.5
$ws.Range("I79").Value = 6699.75
$ws.Range("J79").Value = 7799.8
$ws.Range("K79").Value = 6699.75
$ws.Range("L79").Value = 7799.8
$ws.Range("M79").Value = -5607.75
$ws.Range("N79").Value = -9983.799999999999

$ws.Range("H132").Value = 1360.921
$ws.Range("I132").Value = 1272.5
$ws.Range("K132").Value = 3817.5
$ws.Range("M132").Value = -1287.5

$ws.Range("H137").Value = 3017.36
$ws.Range("I137").Value = 1360.2222
$ws.Range("J137").Value = 3949.5
$ws.Range("K137").Value = 4080.6666
$ws.Range("L137").Value = 11848.5
$ws.Range("M137").Value = -1530.6666
$ws.Range("N137").Value = -16948.5

$ws.Range("H138").Value = 5273.25
$ws.Range("I138").Value = 4731.3335
$ws.Range("K138").Value = 14194.0005
$ws.Range("M138").Value = -9054.000499999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1284.2858
$ws.Range("I2").Value = 1178
$ws.Range("K2").Value = 1178
$ws.Range("M2").Value = -1065

$ws.Range("H32").Value = 16924.305
$ws.Range("I32").Value = 7944.3213
$ws.Range("K32").Value = 7944.3213
$ws.Range("M32").Value = -7657.3213

$ws.Range("H63").Value = 8420.6
$ws.Range("I63").Value = 12250
$ws.Range("J63").Value = 7463.25
$ws.Range("K63").Value = 12250
$ws.Range("L63").Value = 7463.25
$ws.Range("M63").Value = -11564
$ws.Range("N63").Value = -8835.25

$ws.Range("H66").Value = 8420.6
$ws.Range("I66").Value = 12250
$ws.Range("J66").Value = 7463.25
$ws.Range("K66").Value = 61250
$ws.Range("L66").Value = 37316.25
$ws.Range("M66").Value = -57818
$ws.Range("N66").Value = -44180.25

$ws.Range("H88").Value = 1651.0834
$ws.Range("I88").Value = 1587.4
$ws.Range("J88").Value = 1696.5714
$ws.Range("K88").Value = 1587.4
$ws.Range("L88").Value = 1696.5714
$ws.Range("M88").Value = -1181.4
$ws.Range("N88").Value = -2508.5714

$ws.Range("H91").Value = 1651.0834
$ws.Range("I91").Value = 1587.4
$ws.Range("J91").Value = 1696.5714
$ws.Range("K91").Value = 1587.4
$ws.Range("L91").Value = 1696.5714
$ws.Range("M91").Value = -183.4000000000001
$ws.Range("N91").Value = -4504.5714

$ws.Range("H101").Value = 59000
$ws.Range("J101").Value = 59000
$ws.Range("L101").Value = 59000
$ws.Range("N101").Value = -65490

$ws.Range("H116").Value = 1284.2858
$ws.Range("I116").Value = 1178
$ws.Range("K116").Value = 1178
$ws.Range("M116").Value = 1116

$ws.Range("H123").Value = 90000
$ws.Range("J123").Value = 90000
$ws.Range("L123").Value = 90000
$ws.Range("N123").Value = -99800

$ws.Range("H132").Value = 1743.2407
$ws.Range("I132").Value = 1688.88
$ws.Range("K132").Value = 5066.64
$ws.Range("M132").Value = -2536.64

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1284.2858
$ws.Range("I3").Value = 1178
$ws.Range("K3").Value = 1178
$ws.Range("M3").Value = -1064

$ws.Range("H105").Value = 4391.7
$ws.Range("I105").Value = 3872
$ws.Range("J105").Value = 4985.643
$ws.Range("K105").Value = 3872
$ws.Range("L105").Value = 4985.643
$ws.Range("M105").Value = -2125
$ws.Range("N105").Value = -8479.643

$ws.Range("H134").Value = 914
$ws.Range("I134").Value = 738.65515
$ws.Range("K134").Value = 2215.96545
$ws.Range("M134").Value = 319.0345499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4868.1177
$ws.Range("I31").Value = 2544.3333
$ws.Range("K31").Value = 2544.3333
$ws.Range("M31").Value = -2249.3333

$ws.Range("H34").Value = 4868.1177
$ws.Range("I34").Value = 2544.3333
$ws.Range("K34").Value = 2544.3333
$ws.Range("M34").Value = -2342.3333

$ws.Range("H107").Value = 414.5
$ws.Range("I107").Value = 299.92856
$ws.Range("J107").Value = 681.8333
$ws.Range("K107").Value = 299.92856
$ws.Range("L107").Value = 681.8333
$ws.Range("M107").Value = 1620.07144
$ws.Range("N107").Value = -4521.8333

$ws.Range("H132").Value = 2792.75
$ws.Range("I132").Value = 2443.818
$ws.Range("K132").Value = 7331.454000000001
$ws.Range("M132").Value = -4801.454000000001

$ws.Range("H134").Value = 1817.1082
$ws.Range("I134").Value = 1255.3103
$ws.Range("J134").Value = 3853.625
$ws.Range("K134").Value = 3765.9309
$ws.Range("L134").Value = 11560.875
$ws.Range("M134").Value = -1230.9309
$ws.Range("N134").Value = -16630.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1694064.4
$ws.Range("I4").Value = 2445540.8
$ws.Range("J4").Value = 3242.75
$ws.Range("K4").Value = 7336622.399999999
$ws.Range("L4").Value = 9728.25
$ws.Range("M4").Value = -7336510.399999999
$ws.Range("N4").Value = -9952.25

$ws.Range("H15").Value = 584.0833
$ws.Range("I15").Value = 201.5
$ws.Range("J15").Value = 966.6667
$ws.Range("K15").Value = 604.5
$ws.Range("L15").Value = 2900.0001
$ws.Range("M15").Value = -464.5
$ws.Range("N15").Value = -3180.0001

$ws.Range("H75").Value = 739.5
$ws.Range("I75").Value = 549
$ws.Range("J75").Value = 834.75
$ws.Range("K75").Value = 1647
$ws.Range("L75").Value = 2504.25
$ws.Range("M75").Value = -649
$ws.Range("N75").Value = -4500.25

$ws.Range("H78").Value = 739.5
$ws.Range("I78").Value = 549
$ws.Range("J78").Value = 834.75
$ws.Range("K78").Value = 4941
$ws.Range("L78").Value = 7512.75
$ws.Range("M78").Value = 51
$ws.Range("N78").Value = -17496.75

$ws.Range("H114").Value = 300
$ws.Range("J114").Value = 191.44444
$ws.Range("L114").Value = 574.33332
$ws.Range("N114").Value = -7082.33332

$ws.Range("H116").Value = 2530.5
$ws.Range("I116").Value = 29
$ws.Range("K116").Value = 87
$ws.Range("M116").Value = 3355

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9696.571
$ws.Range("I80").Value = 9470.75
$ws.Range("K80").Value = 9470.75
$ws.Range("M80").Value = -8472.75

$ws.Range("H83").Value = 9696.571
$ws.Range("I83").Value = 9470.75
$ws.Range("K83").Value = 47353.75
$ws.Range("M83").Value = -42361.75

$ws.Range("H125").Value = 79999.5
$ws.Range("J125").Value = 79999.5
$ws.Range("L125").Value = 79999.5
$ws.Range("N125").Value = -84919.5

$ws.Range("H132").Value = 2625.6
$ws.Range("I132").Value = 2174.4
$ws.Range("K132").Value = 6523.200000000001
$ws.Range("M132").Value = -3993.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 848.25
$ws.Range("I22").Value = 1000.5
$ws.Range("J22").Value = 696
$ws.Range("K22").Value = 1000.5
$ws.Range("L22").Value = 696
$ws.Range("M22").Value = -705.5
$ws.Range("N22").Value = -1286

$ws.Range("H27").Value = 848.25
$ws.Range("I27").Value = 1000.5
$ws.Range("J27").Value = 696
$ws.Range("K27").Value = 1000.5
$ws.Range("L27").Value = 696
$ws.Range("M27").Value = -893.5
$ws.Range("N27").Value = -910

$ws.Range("H40").Value = 2749.625
$ws.Range("J40").Value = 2997.3333
$ws.Range("L40").Value = 2997.3333
$ws.Range("N40").Value = -3269.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 76749.5
$ws.Range("J135").Value = 76749.5
$ws.Range("L135").Value = 76749.5
$ws.Range("N135").Value = -86889.5

$ws.Range("H136").Value = 25417.953
$ws.Range("J136").Value = 64079.75
$ws.Range("L136").Value = 192239.25
$ws.Range("N136").Value = -197339.25
